{"js": "// Split the combined \"#7cc867#fb5b89#f9cd59#c885da\" highlights paragraph\n// into one paragraph per highlight color, each annotated with its count\n// (bug fix for highlight counts in export files).\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst target = \"#7cc867#fb5b89#f9cd59#c885da\";\nconst replacements = [\n  \"#7cc867: 42\",\n  \"#fb5b89: 18\",\n  \"#f9cd59: 8\",\n  \"#c885da: 8\",\n];\n\nconst hit = paragraphs.items.find((p) => p.text === target);\nif (!hit) {\n  throw new Error(\"Could not find the highlights paragraph to split\");\n}\n\n// Put the first replacement text in the existing paragraph...\nhit.clear();\nhit.insertText(replacements[0], \"Start\");\n\n// ...then insert the remaining ones as new paragraphs right after it,\n// preserving original order.\nlet anchor = hit;\nfor (let i = 1; i < replacements.length; i++) {\n  anchor = anchor.insertParagraph(replacements[i], \"After\");\n}\n\nawait context.sync();\n", "ps1": "# Split the combined \"#7cc867#fb5b89#f9cd59#c885da\" highlights paragraph\n# into one paragraph per highlight color, each annotated with its count\n# (bug fix for highlight counts in export files).\n$d = $word.ActiveDocument\n\n$target = \"#7cc867#fb5b89#f9cd59#c885da\"\n$replacements = @(\n    \"#7cc867: 42\",\n    \"#fb5b89: 18\",\n    \"#f9cd59: 8\",\n    \"#c885da: 8\"\n)\n\n$hit = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text.TrimEnd(\"`r\", \"`a\") -eq $target) {\n        $hit = $p\n        break\n    }\n}\n\nif ($hit -eq $null) {\n    throw \"Could not find the highlights paragraph to split\"\n}\n\n# Overwrite the existing paragraph's text with the first replacement...\n$hit.Range.Text = $replacements[0]\n\n# ...then insert the remaining ones as new paragraphs right after it,\n# preserving original order.\n$cur = $hit\nfor ($i = 1; $i -lt $replacements.Length; $i++) {\n    $cur.Range.InsertParagraphAfter()\n    $cur = $cur.Next()\n    $cur.Range.Text = $replacements[$i]\n}\n"}
